# Apply the "working on other sheets...getting close to being done" edit.
#
# Summary of changes:
#  1. "Data inputs" sheet: new column G ("format") filled in for rows 3-49,
#     new column width for G, sheet becomes the active tab/selected sheet
#     with selection on J20.
#  2. "Data constants" sheet: no longer the active/top-left-pinned sheet;
#     selection moves to J32.
#  3. workbook.xml activeTab follows from (1)/(2) automatically.
#  4. sharedStrings.xml gains a new unique string "general" (used by the
#     new column G values) - this happens automatically as a side effect
#     of setting cell values that aren't already in the shared string table.

$wb = $excel.ActiveWorkbook
$ws5 = $wb.Worksheets.Item("Data inputs")
$ws6 = $wb.Worksheets.Item("Data constants")

# --- New "format" column (G) values for Data inputs, rows 3-49 ---
$formats = @{
    3  = "decimal"
    4  = "general"
    5  = "decimal"
    6  = "decimal"
    7  = "decimal"
    8  = "percentage"
    9  = "percentage"
    10 = "general"
    11 = "general"
    12 = "percentage"
    13 = "general"
    14 = "number"
    15 = "percentage"
    16 = "general"
    17 = "general"
    18 = "general"
    19 = "percentage"
    20 = "general"
    21 = "general"
    22 = "general"
    23 = "percentage"
    24 = "percentage"
    25 = "percentage"
    26 = "percentage"
    27 = "percentage"
    28 = "general"
    29 = "general"
    30 = "percentage"
    31 = "percentage"
    32 = "general"
    33 = "general"
    34 = "general"
    35 = "general"
    36 = "percentage"
    37 = "percentage"
    38 = "percentage"
    39 = "percentage"
    40 = "general"
    41 = "percentage"
    42 = "general"
    43 = "general"
    44 = "general"
    45 = "general"
    46 = "general"
    47 = "general"
    48 = "general"
    49 = "general"
}

foreach ($row in 3..49) {
    $ws5.Range("G$row").Value = $formats[$row]
}

# Widen the new column G on the Data inputs sheet. Target stored width is
# 15.85546875 character-units; this engine snaps ColumnWidth assignments to
# the nearest 1/6 of a character (pixel-grid emulation), so the closest
# achievable stored width is 15 + 5/6 = 15.8333...; 15 lands safely in the
# middle of that snap bucket.
$ws5.Columns.Item(7).ColumnWidth = 15

# Move the selection on "Data constants" before switching away from it so the
# new active sheet ends up being "Data inputs" (matches activeTab going from
# Data constants -> Data inputs), and its own selection updates to J32.
$ws6.Range("J32").Select()

# Activate "Data inputs" as the selected/active sheet with selection J20.
$ws5.Activate()
$ws5.Range("J20").Select()
